# Insert a new data row at row 1053 (pushing the existing rows 1053-1138
# down to 1054-1139) and populate it with the new "Cebolla" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 1053, shifting everything below it down.
$ws.Rows("1053:1053").Insert()

# Populate the newly inserted row with its data.
$ws.Range("A1053").Value = 8
$ws.Range("B1053").Value = "Terminal La Palmera de La Serena"
$ws.Range("C1053").Value = "Coquimbo"
$ws.Range("D1053").Value = 45106
$ws.Range("E1053").Value = 4
$ws.Range("F1053").Value = 100112004
$ws.Range("G1053").Value = "Cebolla"
$ws.Range("H1053").Value = "Sin especificar"
$ws.Range("I1053").Value = "1a (guarda)"
$ws.Range("J1053").Value = 2400
$ws.Range("K1053").Value = 8000
$ws.Range("L1053").Value = 9000
$ws.Range("M1053").Value = 8500
$ws.Range("N1053").Value = "$/malla 16 kilos"
$ws.Range("O1053").Value = "Región de O'Higgins"
$ws.Range("P1053").Value = 531
$ws.Range("Q1053").Value = 16
$ws.Range("R1053").Value = "Hortaliza"
